# "added checks for no description, invalid header"
#
# Adds a new "description" column (header in M1) to the "params" sheet.
# This grows the sheet's used range from A1:L3 to A1:M3, introduces a new
# shared string "description", and leaves the new header cell M1 unstyled
# (matching the rest of the plain header cells such as A1/B1/.../G1).
# Finally, the active selection moves onto the newly-added header cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("params")

$ws.Range("M1").Value = "description"

$null = $ws.Range("M1").Select()
